$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Insert a new row before row 383, shifting existing rows 383:498 down to 384:499
$ws.Rows.Item(383).Insert()

# Populate the new row 383 with data
$ws.Cells.Item(383, 1).Value = 3
$ws.Cells.Item(383, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(383, 3).Value = "Coquimbo"
$ws.Cells.Item(383, 4).Value = 44588
$ws.Cells.Item(383, 5).Value = 5
$ws.Cells.Item(383, 6).Value = "Fruta"
$ws.Cells.Item(383, 7).Value = 100108
$ws.Cells.Item(383, 8).Value = "Tropicales y subtropicales"
$ws.Cells.Item(383, 9).Value = 100108005
$ws.Cells.Item(383, 10).Value = "Piña"
$ws.Cells.Item(383, 11).Value = "Caramelo"
$ws.Cells.Item(383, 12).Value = "Primera"
$ws.Cells.Item(383, 13).Value = 230
$ws.Cells.Item(383, 14).Value = 17500
$ws.Cells.Item(383, 15).Value = 18000
$ws.Cells.Item(383, 16).Value = 17739
$ws.Cells.Item(383, 17).Value = "$/caja 12 unidades"
$ws.Cells.Item(383, 18).Value = "Ecuador"
$ws.Cells.Item(383, 19).Value = 1478
$ws.Cells.Item(383, 20).Value = 12
